# Apply updated crypto price/volume data per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Prefix with an apostrophe so Excel stores the literal text instead of
    # auto-converting number-like strings (e.g. '314.21') to a numeric value,
    # then strip the formatting nudge this causes so the cell keeps its
    # original (default) style.
    $range.Value2 = "'" + $text
    $range.ClearFormats()
}

Set-TextCell $ws.Range("D2") "27.294.75"
Set-TextCell $ws.Range("E2") "  +1.21%  "

Set-TextCell $ws.Range("D3") "1.854.76"
Set-TextCell $ws.Range("E3") "  +1.54%  "

Set-TextCell $ws.Range("E4") "  -0.58%  "

Set-TextCell $ws.Range("D5") "314.21"
Set-TextCell $ws.Range("E5") "  +0.56%  "

Set-TextCell $ws.Range("D6") "1.001"
Set-TextCell $ws.Range("E6") "  -0.52%  "

Set-TextCell $ws.Range("D7") "0.4604"
Set-TextCell $ws.Range("E7") "  +0.46%  "

Set-TextCell $ws.Range("E8") "  +0.23%  "

Set-TextCell $ws.Range("D9") "0.07299"
Set-TextCell $ws.Range("E9") "  -0.12%  "

Set-TextCell $ws.Range("D10") "0.8879"
Set-TextCell $ws.Range("E10") "  +1.76%  "

Set-TextCell $ws.Range("D11") "20.15"
Set-TextCell $ws.Range("E11") "  +2.04%  "

Set-TextCell $ws.Range("D12") "0.07822"
Set-TextCell $ws.Range("E12") "  -1.74%  "

Set-TextCell $ws.Range("D13") "5.388"
Set-TextCell $ws.Range("E13") "  +1.14%  "

Set-TextCell $ws.Range("D14") "1.790.04"
Set-TextCell $ws.Range("E14") "  -3.91%  "

Set-TextCell $ws.Range("D15") "6.532"
Set-TextCell $ws.Range("E15") "  -0.61%  "

Set-TextCell $ws.Range("D16") "91.44"
Set-TextCell $ws.Range("E16") "  +0.17%  "

Set-TextCell $ws.Range("E17") "  -0.54%  "

Set-TextCell $ws.Range("D18") "0.000008933"
Set-TextCell $ws.Range("E18") "  +0.90%  "

Set-TextCell $ws.Range("E19") "  -0.68%  "

Set-TextCell $ws.Range("D20") "14.77"
Set-TextCell $ws.Range("E20") "  +0.37%  "

Set-TextCell $ws.Range("D21") "27.305.50"
Set-TextCell $ws.Range("E21") "  -0.09%  "

Set-TextCell $ws.Range("D22") "5.113"

Set-TextCell $ws.Range("D23") "10.56"
Set-TextCell $ws.Range("E23") "  +0.23%  "

Set-TextCell $ws.Range("D24") "2.050.49"
Set-TextCell $ws.Range("E24") "  -4.49%  "

Set-TextCell $ws.Range("D25") "1.924"
Set-TextCell $ws.Range("E25") "  +4.03%  "

Set-TextCell $ws.Range("D26") "152.09"
Set-TextCell $ws.Range("E26") "  -0.60%  "

Set-TextCell $ws.Range("D27") "18.48"
Set-TextCell $ws.Range("E27") "  +0.58%  "

Set-TextCell $ws.Range("D28") "2.062"
Set-TextCell $ws.Range("E28") "  +0.74%  "

Set-TextCell $ws.Range("E29") "  +0.76%  "

Set-TextCell $ws.Range("D30") "5.072"
Set-TextCell $ws.Range("E30") "  -1.37%  "

Set-TextCell $ws.Range("D31") "0.08826"
Set-TextCell $ws.Range("E31") "  -0.40%  "

Set-TextCell $ws.Range("B32") "ImmutableX"
Set-TextCell $ws.Range("C32") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws.Range("D32") "0.7737"
Set-TextCell $ws.Range("E32") "  +6.17%  "

Set-TextCell $ws.Range("B33") "HuobiToken"
Set-TextCell $ws.Range("C33") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws.Range("D33") "3.077"
Set-TextCell $ws.Range("E33") "  +4.05%  "

Set-TextCell $ws.Range("D34") "1.171"
Set-TextCell $ws.Range("E34") "  +3.75%  "

Set-TextCell $ws.Range("D35") "4.507"
Set-TextCell $ws.Range("E35") "  +1.46%  "

Set-TextCell $ws.Range("D36") "2.760"
Set-TextCell $ws.Range("E36") "  +13.33%  "

Set-TextCell $ws.Range("D37") "1.081"
Set-TextCell $ws.Range("E37") "  +0.87%  "

Set-TextCell $ws.Range("E38") "  +0.67%  "

Set-TextCell $ws.Range("D39") "0.05252"
Set-TextCell $ws.Range("E39") "  +0.49%  "

Set-TextCell $ws.Range("D40") "2.953"
Set-TextCell $ws.Range("E40") "  +0.59%  "

Set-TextCell $ws.Range("D41") "7.065"
Set-TextCell $ws.Range("E41") "  -1.42%  "

Set-TextCell $ws.Range("D42") "0.5124"
Set-TextCell $ws.Range("E42") "  -0.14%  "

Set-TextCell $ws.Range("D43") "0.1635"
Set-TextCell $ws.Range("E43") "  +0.49%  "

Set-TextCell $ws.Range("D44") "8.409"
Set-TextCell $ws.Range("E44") "  +2.69%  "

Set-TextCell $ws.Range("D45") "0.4799"
Set-TextCell $ws.Range("E45") "  -0.56%  "

Set-TextCell $ws.Range("D46") "10.34"
Set-TextCell $ws.Range("E46") "  +1.72%  "

Set-TextCell $ws.Range("D47") "1.001"
Set-TextCell $ws.Range("E47") "  -0.59%  "

Set-TextCell $ws.Range("D48") "102.48"
Set-TextCell $ws.Range("E48") "  +0.00%  "

Set-TextCell $ws.Range("D49") "1.644"
Set-TextCell $ws.Range("E49") "  +0.95%  "

Set-TextCell $ws.Range("E50") "  +0.14%  "

Set-TextCell $ws.Range("D51") "65.82"
Set-TextCell $ws.Range("E51") "  +1.40%  "
